$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
$d.Name = "Hello"
Write-Output "renamed"
